# Weekly reports and Timesheets
# Update the "Weekly" timesheet from the week of 16/02 to the week of 11/05,
# and correct the "Week of:" date to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day-of-week labels for the new week (Sun 11/05 .. Sat 17/05)
$ws.Range("A11").Value = "Sun 11/05"
$ws.Range("A12").Value = "Mon 12/05"
$ws.Range("A13").Value = "Tue 13/05"
$ws.Range("A14").Value = "Wed 14 /05"
$ws.Range("A15").Value = "Thur 15 /05"
$ws.Range("A16").Value = "Fri   16/05"
$ws.Range("A17").Value = "Sat 17/05"

# "Week of:" date corrected to the Sunday that starts the new week (11 May 2014)
$ws.Range("G8").Value = 41770

# Leave the cursor/selection where the edit finished, on the corrected date cell
$ws.Range("G8:H8").Select() | Out-Null
